# Daily update at 8 AM UTC
# Adds the next day's row of data to the "Wins Over Time" log and moves the
# "most recent row" date formatting down to the newly-appended row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previously-last row (67) used a distinct date-only number format to
# highlight it as the latest entry. Since a new row is being appended below
# it, restore row 67's date column to the regular format used by all the
# other historical rows.
$ws.Range("A67").NumberFormat = $ws.Range("A66").NumberFormat

# Append the new day's data as row 68.
$ws.Range("A68").Value = 45808
$ws.Range("B68").Value = 291
$ws.Range("C68").Value = 287
$ws.Range("D68").Value = 291

# Give the new last row the "latest entry" date-only number format that row
# 67 previously had.
$ws.Range("A68").NumberFormat = "YYYY-MM-DD"
